$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 / column E ("time") held a raw date-time serial number (43530.5902777778,
# formatted as MM/DD/YYYY HH:MM:SS). The uploaded workbook instead stores that cell
# as plain text reading "27/5/2019 15:55" -- the same string already used by the
# cell directly above it (E2) -- with ordinary General formatting.
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E3").Value = "27/5/2019 15:55"

# The previously selected/active cell in the sheet view was E4; the new workbook
# shows E3 as the active selection instead.
$ws.Range("E3").Select()
